$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Add a new row at the end of the table for the new video entry.
$newRow = $table.Rows.Add()
$rowIndex = $table.Rows.Count

# Column 1: description of the video.
$cell1 = $table.Cell($rowIndex, 1)
$cell1.Range.Text = "Changing your data structure to store values in separate columns by categories that are currently stored in a column then using those data in to build a confidence interval. "

# Column 2: hyperlink to the video, followed by a trailing space
# (matching the pattern used by the other rows in the table).
$cell2 = $table.Cell($rowIndex, 2)
$linkRange = $cell2.Range
$linkRange.End = $linkRange.End - 1
$linkRange.Text = "https://youtu.be/Si1Q5h1w-Xo"
$d.Hyperlinks.Add($linkRange, "https://youtu.be/Si1Q5h1w-Xo", [Type]::Missing, [Type]::Missing, "https://youtu.be/Si1Q5h1w-Xo")
$cell2b = $table.Cell($rowIndex, 2)
$cell2b.Range.InsertAfter(" ")

# Column 3: creation date.
$cell3 = $table.Cell($rowIndex, 3)
$cell3.Range.Text = "Created 4/3/24"

# Column 4: associated stone references.
$cell4 = $table.Cell($rowIndex, 4)
$cell4.Range.Text = "Stones 21, 47 "
